$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix the B2 header label (was "unnamed: 1_level_1", should match B1 "total")
$ws.Range("B2").Value = $ws.Range("B1").Value2

# Remove the "situação do domicílio" sub-header row (row 5)
$ws.Rows.Item(5).Delete()

# Remove the "grandes regiões e unidades da federação" sub-header row (now row 7 after the previous delete)
$ws.Rows.Item(7).Delete()
